$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$DEFAULT = -17.57101932020553

$rows = @(
    @($DEFAULT, 2.040425197981548, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 3.01050585659935, $DEFAULT, $DEFAULT),
    @($DEFAULT, 2.322692520630028, 2.167872843743806, $DEFAULT, 3.399786627464748, $DEFAULT, $DEFAULT, $DEFAULT, 2.559673147834696, $DEFAULT),
    @($DEFAULT, 1.823964545617049, $DEFAULT, $DEFAULT, $DEFAULT, 2.903788740627212, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @(2.634674910386622, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, 1.809579138673929, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @(3.785504113268936, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 1.191586564924966, $DEFAULT, 1.756754780915831),
    @($DEFAULT, $DEFAULT, $DEFAULT, 2.839387517811856, $DEFAULT, 2.52970589102887, $DEFAULT, $DEFAULT, $DEFAULT, 1.753110672823682),
    @($DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, 2.406077492938622, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 2.013666147195698, 1.943511750656977),
    @($DEFAULT, $DEFAULT, 1.330875781841625, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 2.150455410167084),
    @($DEFAULT, $DEFAULT, 1.263590605206954, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 2.405602136221352, $DEFAULT),
    @($DEFAULT, 1.861353341319296, 2.186678315141741, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 1.302199473776125, 1.338279040124827, $DEFAULT),
    @($DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 0.6662326612693678, 1.162886418125556, $DEFAULT),
    @($DEFAULT, $DEFAULT, 1.594439762853831, $DEFAULT, $DEFAULT, $DEFAULT, $DEFAULT, 1.668143310426421, $DEFAULT, $DEFAULT),
    @($DEFAULT, 0.7797688247501143, 1.593359792852021, $DEFAULT, 3.239612716336004, $DEFAULT, $DEFAULT, 1.278816399059651, $DEFAULT, 2.31139237936753),
    @($DEFAULT, 0.9947558885406091, $DEFAULT, 2.013299133734759, $DEFAULT, 2.753065612696394, 4.321921056140708, $DEFAULT, $DEFAULT, $DEFAULT),
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowvals = $rows[$i]
    for ($j = 0; $j -lt $rowvals.Count; $j++) {
        $ws.Cells.Item($r, $j + 2).Value = $rowvals[$j]
    }
}